# Auto-generated Word COM-interop script.
# Applies proofing-pass run splits (w:proofErr spellStart/spellEnd,
# gramStart/gramEnd) and highlight formatting fixes to match the target revision.
$d = $word.ActiveDocument

$xml4 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="green"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>network partner (default to RUMC)</w:t></w:r><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve"> -  </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>datalist</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$d.Paragraphs(4).Range.InsertXML($xml4)

$xml5 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="green"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>Agency (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>defualt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve"> to 1039a)</w:t></w:r><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve"> - </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>datalist</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$d.Paragraphs(5).Range.InsertXML($xml5)

$xml6 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="green"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>Donor Name (e.g. Jewel, Sam''s Club)</w:t></w:r><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve"> – </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>datalist</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$d.Paragraphs(6).Range.InsertXML($xml6)

$xml7 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="green"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>City</w:t></w:r><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve"> - </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>datalist</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$d.Paragraphs(7).Range.InsertXML($xml7)

$xml19 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="720"/><w:rPr><w:highlight w:val="green"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve">Validate those </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve">donations </w:t></w:r><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve"> in</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve"> the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>js</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$d.Paragraphs(19).Range.InsertXML($xml19)

$xml22 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve">Submit donations in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>php</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$d.Paragraphs(22).Range.InsertXML($xml22)

$xml24 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:highlight w:val="green"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve">Name - </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>datalist</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$d.Paragraphs(24).Range.InsertXML($xml24)

$xml25 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:highlight w:val="green"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve">City - </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>datalist</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$d.Paragraphs(25).Range.InsertXML($xml25)

$xml27 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:highlight w:val="green"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve">Zip - </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>datalist</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$d.Paragraphs(27).Range.InsertXML($xml27)

$xml28 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:highlight w:val="green"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve">Address - </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>datalist</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$d.Paragraphs(28).Range.InsertXML($xml28)

$xml29 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:highlight w:val="green"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>phoneNumber</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve"> – </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>datalist</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$d.Paragraphs(29).Range.InsertXML($xml29)

$xml30 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/></w:r><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve">Validate </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>donpar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve"> in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>js</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$d.Paragraphs(30).Range.InsertXML($xml30)

$xml31 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/></w:r><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve">Submit </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>donpar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve"> in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>php</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$d.Paragraphs(31).Range.InsertXML($xml31)

$xml37 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:highlight w:val="green"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>Delete donations</w:t></w:r></w:p>'
$d.Paragraphs(37).Range.InsertXML($xml37)

$xml38 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>Delete donation partners</w:t></w:r></w:p>'
$d.Paragraphs(38).Range.InsertXML($xml38)

